# Applies updated odds values to Sheet1 per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 8).Value = 6.8
$ws.Cells.Item(2, 18).Value = 1.57
$ws.Cells.Item(2, 19).Value = 2.64
$ws.Cells.Item(2, 21).Value = 2.1
$ws.Cells.Item(2, 38).Value = 32
$ws.Cells.Item(3, 6).Value = 2.02
$ws.Cells.Item(3, 7).Value = 2.06
$ws.Cells.Item(3, 8).Value = 3.65
$ws.Cells.Item(3, 9).Value = 3.8
$ws.Cells.Item(3, 13).Value = 1.05
$ws.Cells.Item(3, 14).Value = 4.9
$ws.Cells.Item(3, 21).Value = 2.46
$ws.Cells.Item(3, 24).Value = 21
$ws.Cells.Item(3, 25).Value = 18
$ws.Cells.Item(3, 26).Value = 1000
$ws.Cells.Item(3, 27).Value = 70
$ws.Cells.Item(3, 30).Value = 15.5
$ws.Cells.Item(3, 31).Value = 48
$ws.Cells.Item(3, 32).Value = 15
$ws.Cells.Item(3, 33).Value = 11
$ws.Cells.Item(3, 35).Value = 980
$ws.Cells.Item(3, 37).Value = 25
$ws.Cells.Item(3, 38).Value = 48
$ws.Cells.Item(3, 40).Value = 11
$ws.Cells.Item(3, 41).Value = 32
$ws.Cells.Item(4, 6).Value = 2.44
$ws.Cells.Item(4, 7).Value = 3.4
$ws.Cells.Item(4, 8).Value = 2.3
$ws.Cells.Item(4, 10).Value = 3.35
$ws.Cells.Item(4, 11).Value = 6.8
$ws.Cells.Item(5, 6).Value = 1.48
$ws.Cells.Item(5, 7).Value = 1.6
$ws.Cells.Item(5, 8).Value = 1.09
$ws.Cells.Item(5, 9).Value = 9.199999999999999
$ws.Cells.Item(5, 10).Value = 3.9
$ws.Cells.Item(5, 11).Value = 5.1
$ws.Cells.Item(6, 14).Value = 3.5
$ws.Cells.Item(6, 16).Value = 1.82
$ws.Cells.Item(6, 17).Value = 2.2
$ws.Cells.Item(6, 26).Value = 23
$ws.Cells.Item(6, 28).Value = 9.6
$ws.Cells.Item(6, 31).Value = 100
$ws.Cells.Item(6, 36).Value = 34
$ws.Cells.Item(6, 38).Value = 1000
$ws.Cells.Item(6, 40).Value = 23
$ws.Cells.Item(6, 41).Value = 44
$ws.Cells.Item(7, 6).Value = 2.54
$ws.Cells.Item(7, 17).Value = 1.54
$ws.Cells.Item(7, 18).Value = 1.72
$ws.Cells.Item(7, 19).Value = 2.34
$ws.Cells.Item(7, 41).Value = 14
$ws.Cells.Item(8, 16).Value = 2.52
$ws.Cells.Item(8, 19).Value = 2.54
$ws.Cells.Item(8, 21).Value = 1.96
$ws.Cells.Item(9, 14).Value = 6.4
$ws.Cells.Item(9, 18).Value = 1.76
$ws.Cells.Item(9, 29).Value = 19
$ws.Cells.Item(9, 35).Value = 300
$ws.Cells.Item(10, 7).Value = 1.31
$ws.Cells.Item(10, 8).Value = 11
$ws.Cells.Item(10, 10).Value = 6.6
$ws.Cells.Item(10, 11).Value = 7
$ws.Cells.Item(10, 14).Value = 8.800000000000001
$ws.Cells.Item(10, 15).Value = 1.12
$ws.Cells.Item(10, 17).Value = 1.37
$ws.Cells.Item(10, 18).Value = 2
$ws.Cells.Item(10, 19).Value = 1.92
$ws.Cells.Item(10, 21).Value = 2.26
$ws.Cells.Item(10, 24).Value = 46
$ws.Cells.Item(10, 25).Value = 60
$ws.Cells.Item(10, 28).Value = 15.5
$ws.Cells.Item(10, 29).Value = 17
$ws.Cells.Item(10, 30).Value = 100
$ws.Cells.Item(10, 31).Value = 170
$ws.Cells.Item(10, 35).Value = 1000
$ws.Cells.Item(10, 37).Value = 13.5
$ws.Cells.Item(10, 38).Value = 29
$ws.Cells.Item(10, 39).Value = 1000
$ws.Cells.Item(10, 41).Value = 140
$ws.Cells.Item(11, 6).Value = 6
$ws.Cells.Item(11, 10).Value = 4.2
$ws.Cells.Item(11, 13).Value = 1.06
$ws.Cells.Item(11, 14).Value = 4.3
$ws.Cells.Item(11, 16).Value = 2.12
$ws.Cells.Item(11, 17).Value = 1.84
$ws.Cells.Item(11, 18).Value = 1.45
$ws.Cells.Item(11, 20).Value = 1.89
$ws.Cells.Item(11, 25).Value = 9.4
$ws.Cells.Item(11, 32).Value = 160
$ws.Cells.Item(11, 35).Value = 36
$ws.Cells.Item(11, 36).Value = 210
$ws.Cells.Item(11, 38).Value = 80
$ws.Cells.Item(11, 40).Value = 120
$ws.Cells.Item(12, 8).Value = 2.48
$ws.Cells.Item(12, 14).Value = 4.7
$ws.Cells.Item(12, 16).Value = 2.24
$ws.Cells.Item(12, 17).Value = 1.76
$ws.Cells.Item(12, 18).Value = 1.49
$ws.Cells.Item(12, 19).Value = 2.82
$ws.Cells.Item(12, 20).Value = 1.64
$ws.Cells.Item(12, 24).Value = 19
$ws.Cells.Item(12, 25).Value = 13
$ws.Cells.Item(12, 31).Value = 25
$ws.Cells.Item(12, 32).Value = 22
